$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "date-looking" text in column B to stay as literal text instead
# of being auto-converted to a date serial number, matching the rest of the
# sheet (which stores dates as plain strings).
$ws.Range("B3:B4").NumberFormat = "@"

$ws.Range("A3").Value = "222BBB222"
$ws.Range("B3").Value = "2025-10-21"
$ws.Range("C3").Value = "Anna Nagar"
$ws.Range("D3").Value = "BBB"

$ws.Range("A4").Value = "333CCC333"
$ws.Range("B4").Value = "2025-10-22"
$ws.Range("C4").Value = "Perungudi"
$ws.Range("D4").Value = "CCC"
